# Update last_edited_time (column D) values on the "Lũy kế tháng CẦN THƠ" sheet.
#
# - Rows 4, 5, 6, 8, 12, 13: last_edited_time changes from
#   "2024-08-30T20:17:00.000Z" -> "2024-08-31T05:43:00.000Z"
# - Row 7 (the "Tháng 8" record): last_edited_time changes from
#   "2024-08-30T20:17:00.000Z" -> "2024-08-31T05:40:00.000Z"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lũy kế tháng CẦN THƠ")

$newValueCommon = "2024-08-31T05:43:00.000Z"
$newValueRow7 = "2024-08-31T05:40:00.000Z"

$rowsCommon = @(4, 5, 6, 8, 12, 13)
foreach ($r in $rowsCommon) {
    $ws.Cells.Item($r, 4).Value = $newValueCommon
}

$ws.Cells.Item(7, 4).Value = $newValueRow7
